# Weekly fruit/vegetable price update: insert two new weekly records
# (rows 256-257) into the Apio sheet, pushing the existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 256, shifting rows 256:269 down to 258:271.
$ws.Rows("256:257").Insert()

# Row 256: new "Primera" quality record for date 44714
$ws.Range("A256").Value = 11
$ws.Range("B256").Value = "Vega Monumental Concepción"
$ws.Range("C256").Value = "Bíobío"
$ws.Range("D256").Value = 44714
$ws.Range("E256").Value = 8
$ws.Range("F256").Value = 100112017
$ws.Range("G256").Value = "Apio"
$ws.Range("H256").Value = "Americana (o)"
$ws.Range("I256").Value = "Primera"
$ws.Range("J256").Value = 350
$ws.Range("K256").Value = 7000
$ws.Range("L256").Value = 7500
$ws.Range("M256").Value = 7214
$ws.Range("N256").Value = "$/docena de matas"
$ws.Range("O256").Value = "Región de Coquimbo"
$ws.Range("P256").Value = 1202
$ws.Range("Q256").Value = 6
$ws.Range("R256").Value = "Hortaliza"

# Row 257: new "Segunda" quality record for date 44714
$ws.Range("A257").Value = 11
$ws.Range("B257").Value = "Vega Monumental Concepción"
$ws.Range("C257").Value = "Bíobío"
$ws.Range("D257").Value = 44714
$ws.Range("E257").Value = 8
$ws.Range("F257").Value = 100112017
$ws.Range("G257").Value = "Apio"
$ws.Range("H257").Value = "Americana (o)"
$ws.Range("I257").Value = "Segunda"
$ws.Range("J257").Value = 220
$ws.Range("K257").Value = 5000
$ws.Range("L257").Value = 5500
$ws.Range("M257").Value = 5273
$ws.Range("N257").Value = "$/docena de matas"
$ws.Range("O257").Value = "Región de Coquimbo"
$ws.Range("P257").Value = 879
$ws.Range("Q257").Value = 6
$ws.Range("R257").Value = "Hortaliza"
